$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new daily rows (256-269), continuing the date series (serial dates
# 44330-44343, i.e. 2021-05-14 through 2021-05-27) with zero counts, matching
# the style/format of the existing rows above (row 255 as template).

$startRow = 256
$startSerial = 44330
$endSerial = 44343

$row = $startRow
for ($serial = $startSerial; $serial -le $endSerial; $serial++) {
    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $row++
}

$lastRow = $row - 1

# Copy the style of the previous data row (255) to the newly added rows so
# formatting (date number format, borders, alignment) matches.
$srcRange = $ws.Range("A255:D255")
$dstRange = $ws.Range("A256:D$lastRow")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122) # xlPasteFormats
